$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.019.25"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.846.96"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").Value = "'309.52"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "'0.4761"
$ws.Range("E7").Value = "  +1.95%  "
$ws.Range("E8").Value = "  +1.53%  "
$ws.Range("D9").Value = "'0.07242"
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("D10").Value = "'0.9300"
$ws.Range("E10").Value = "  +1.68%  "
$ws.Range("D11").Value = "'19.84"
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("D12").Value = "'0.07776"
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D13").Value = "1.900.34"
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("D14").Value = "'5.392"
$ws.Range("E14").Value = "  +2.12%  "
$ws.Range("D15").Value = "'6.470"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").Value = "'88.94"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "'1.018"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "'0.000008665"
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").Value = "27.016.43"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "'14.55"
$ws.Range("E21").Value = "  +1.34%  "
$ws.Range("D22").Value = "'5.051"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").Value = "'10.65"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "'1.927"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'152.86"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").Value = "'18.26"
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("D27").Value = "'1.986"
$ws.Range("E27").Value = "  -3.05%  "
$ws.Range("D28").Value = "'114.50"
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("D29").Value = "'4.947"
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").Value = "'0.08866"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  +4.01%  "
$ws.Range("D32").Value = "'1.179"
$ws.Range("E32").Value = "  +0.44%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.509"
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7368"
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("D35").Value = "'2.651"
$ws.Range("E36").Value = "  +3.37%  "
$ws.Range("D37").Value = "'0.01969"
$ws.Range("E37").Value = "  +1.51%  "
$ws.Range("D38").Value = "'0.05258"
$ws.Range("E38").Value = "  +1.76%  "
$ws.Range("D39").Value = "'2.977"
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("E40").Value = "  +1.38%  "
$ws.Range("D41").Value = "'7.035"
$ws.Range("E41").Value = "  +1.91%  "
$ws.Range("D42").Value = "'0.1516"
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("D43").Value = "'8.283"
$ws.Range("E43").Value = "  +1.82%  "
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("D47").Value = "'101.58"
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("D48").Value = "'1.608"
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("D49").Value = "'65.63"
$ws.Range("E49").Value = "  +1.23%  "
$ws.Range("D50").Value = "'0.06061"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").Value = "'0.8921"
$ws.Range("E51").Value = "  +3.43%  "
